{"js": "// Replace each two-digit multiplication problem's text with its updated\n// value. Every occurrence in the source document is unique, so searching\n// the whole body for the old text and replacing it is unambiguous.\n\nconst pairs = [\n  [\"60\u00d794=\", \"63\u00d769=\"],\n  [\"69\u00d748=\", \"32\u00d745=\"],\n  [\"52\u00d713=\", \"77\u00d769=\"],\n  [\"46\u00d794=\", \"34\u00d784=\"],\n  [\"40\u00d741=\", \"34\u00d795=\"],\n  [\"66\u00d716=\", \"74\u00d730=\"],\n  [\"44\u00d715=\", \"42\u00d729=\"],\n  [\"63\u00d738=\", \"75\u00d781=\"],\n  [\"17\u00d758=\", \"93\u00d766=\"],\n  [\"73\u00d721=\", \"42\u00d749=\"],\n  [\"95\u00d770=\", \"72\u00d758=\"],\n  [\"12\u00d790=\", \"89\u00d768=\"],\n  [\"95\u00d787=\", \"39\u00d786=\"],\n  [\"50\u00d771=\", \"75\u00d723=\"],\n  [\"80\u00d764=\", \"81\u00d791=\"],\n  [\"28\u00d759=\", \"62\u00d785=\"],\n  [\"56\u00d738=\", \"43\u00d745=\"],\n  [\"41\u00d715=\", \"54\u00d753=\"],\n  [\"87\u00d738=\", \"83\u00d758=\"],\n  [\"36\u00d716=\", \"63\u00d747=\"],\n  [\"74\u00d738=\", \"82\u00d733=\"],\n  [\"83\u00d767=\", \"53\u00d736=\"],\n  [\"84\u00d792=\", \"49\u00d751=\"],\n  [\"42\u00d753=\", \"80\u00d762=\"],\n  [\"14\u00d767=\", \"83\u00d734=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication problem's text with its updated\n# value. Every occurrence in the source document is unique, so a simple\n# whole-document Find/Replace per pair is unambiguous and order-independent.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"60\u00d794=\", \"63\u00d769=\"),\n    @(\"69\u00d748=\", \"32\u00d745=\"),\n    @(\"52\u00d713=\", \"77\u00d769=\"),\n    @(\"46\u00d794=\", \"34\u00d784=\"),\n    @(\"40\u00d741=\", \"34\u00d795=\"),\n    @(\"66\u00d716=\", \"74\u00d730=\"),\n    @(\"44\u00d715=\", \"42\u00d729=\"),\n    @(\"63\u00d738=\", \"75\u00d781=\"),\n    @(\"17\u00d758=\", \"93\u00d766=\"),\n    @(\"73\u00d721=\", \"42\u00d749=\"),\n    @(\"95\u00d770=\", \"72\u00d758=\"),\n    @(\"12\u00d790=\", \"89\u00d768=\"),\n    @(\"95\u00d787=\", \"39\u00d786=\"),\n    @(\"50\u00d771=\", \"75\u00d723=\"),\n    @(\"80\u00d764=\", \"81\u00d791=\"),\n    @(\"28\u00d759=\", \"62\u00d785=\"),\n    @(\"56\u00d738=\", \"43\u00d745=\"),\n    @(\"41\u00d715=\", \"54\u00d753=\"),\n    @(\"87\u00d738=\", \"83\u00d758=\"),\n    @(\"36\u00d716=\", \"63\u00d747=\"),\n    @(\"74\u00d738=\", \"82\u00d733=\"),\n    @(\"83\u00d767=\", \"53\u00d736=\"),\n    @(\"84\u00d792=\", \"49\u00d751=\"),\n    @(\"42\u00d753=\", \"80\u00d762=\"),\n    @(\"14\u00d767=\", \"83\u00d734=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute([ref]$old, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$new, [ref]2) | Out-Null\n}\n"}
